# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh Japon's (row 32) daily numbers
# - Zimbabue's case count overtakes several countries, so it moves up the
#   (cases-desc sorted) table from row 179 to row 172; every row that used
#   to sit between the old and new position shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 05:52"

# --- Row 32: Japon ------------------------------------------------------
$ws.Range("B32").Value = 13965
$ws.Range("C32").Value = 70
$ws.Range("D32").Value = 2368
$ws.Range("E32").Value = 11172
$ws.Range("F32").Value = 306
$ws.Range("G32").Value = 12
$ws.Range("H32").Value = 425

# --- Rows 172-179: Zimbabue moves up the ranking ------------------------
# Row 172 becomes Zimbabue with its freshly updated figures.
$ws.Range("A172").Value = "Zimbabue"
$ws.Range("B172").Value = 40
$ws.Range("C172").Value = 8
$ws.Range("D172").Value = 5
$ws.Range("E172").Value = 31
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 4

# Row 173 becomes Puerto Rico (previously row 172's data).
$ws.Range("A173").Value = "Puerto Rico"
$ws.Range("B173").Value = 39
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 36
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 2

# Row 174 becomes Eritrea (previously row 173's data).
$ws.Range("A174").Value = "Eritrea"
$ws.Range("B174").Value = 39
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 19
$ws.Range("E174").Value = 20
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Row 175 becomes Mongolia (previously row 174's data).
$ws.Range("A175").Value = "Mongolia"
$ws.Range("B175").Value = 38
$ws.Range("C175").Value = 0
$ws.Range("D175").Value = 10
$ws.Range("E175").Value = 28
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Row 176 becomes San Martin (Parte Francesa) (previously row 175's data).
$ws.Range("A176").Value = "San Martin (Parte Francesa)"
$ws.Range("B176").Value = 38
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 24
$ws.Range("E176").Value = 11
$ws.Range("F176").Value = 3
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 3

# Row 177 becomes Malaui (previously row 176's data).
$ws.Range("A177").Value = "Malaui"
$ws.Range("B177").Value = 36
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 7
$ws.Range("E177").Value = 26
$ws.Range("F177").Value = 1
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 3

# Row 178 becomes Sudan del Sur (previously row 177's data).
$ws.Range("A178").Value = "Sudan del Sur"
$ws.Range("B178").Value = 34
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 34
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 179 becomes Guam (previously row 178's data).
$ws.Range("A179").Value = "Guam"
$ws.Range("B179").Value = 32
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 31
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1

# Row 180 (Angola) is unchanged.
